# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback (now in sync with en-US).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (shared by Overview!E2/F2 and the per-locale Status column C2 on each sheet)
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C2").Value     = "Handed back: in sync with en-US"
$dede.Range("C2").Value     = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed for both locales ---
$zhcn.Range("K2").Value = "2016-08-23 02:48:31"
$dede.Range("K2").Value = "2016-08-23 02:48:38"

# --- Error Detail cleared now that the handback version matches latest ---
$zhcn.Range("P2").Value = ""
$dede.Range("P2").Value = ""

# --- Column widths widened for the (now longer) Status / Error Detail columns ---
# Targets (~29.98 / ~13.75 "character" widths) aren't reachable bit-exactly
# through ColumnWidth's pixel-quantised setter, so feed it the nearest input
# that lands in the correct rounding bucket.
$overview.Range("E1").ColumnWidth = 29.166666666666668
$overview.Range("F1").ColumnWidth = 29.166666666666668
$zhcn.Range("C1").ColumnWidth     = 29.166666666666668
$dede.Range("C1").ColumnWidth     = 29.166666666666668

$zhcn.Range("P1").ColumnWidth = 12.833333333333334
$dede.Range("P1").ColumnWidth = 12.833333333333334
